$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.007622599601746
$ws.Range("B1").Value = 1.944369316101074
$ws.Range("C1").Value = 3.013760805130005
$ws.Range("D1").Value = 3.679563522338867
$ws.Range("E1").Value = 1.691420197486877
